$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers (F1, G1, X1, Y1, Z1)
$ws.Range("F1").Value = "FT_Goals_H"
$ws.Range("G1").Value = "FT_Goals_A"
$ws.Range("X1").Value = "FT_Odds_H"
$ws.Range("Y1").Value = "FT_Odds_D"
$ws.Range("Z1").Value = "FT_Odds_A"

# Ensure new Date cells stay as plain text (otherwise Excel auto-parses
# "DD/MM/YYYY"-looking strings into date serials).
foreach ($addr in @("B107","B108","B109","B110","B111","B112","B113")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Add the 7 new match rows (107-113)
# Row 107: OFI Crete vs Volos NFC (03/01/2023)
$ws.Range("A107").Value = "G1"
$ws.Range("B107").Value = "03/01/2023"
$ws.Range("C107").Value = "14:30"
$ws.Range("D107").Value = "OFI Crete"
$ws.Range("E107").Value = "Volos NFC"
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = "D"
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = "D"
$ws.Range("L107").Value = 22
$ws.Range("M107").Value = 12
$ws.Range("N107").Value = 7
$ws.Range("O107").Value = 3
$ws.Range("P107").Value = 16
$ws.Range("Q107").Value = 14
$ws.Range("R107").Value = 8
$ws.Range("S107").Value = 1
$ws.Range("T107").Value = 3
$ws.Range("U107").Value = 2
$ws.Range("V107").Value = 0
$ws.Range("W107").Value = 0
$ws.Range("X107").Value = 2.25
$ws.Range("Y107").Value = 3.3
$ws.Range("Z107").Value = 3.2
$ws.Range("AA107").Value = 2.1
$ws.Range("AB107").Value = 1.7

# Row 108: Levadeiakos vs Panathinaikos (03/01/2023)
$ws.Range("A108").Value = "G1"
$ws.Range("B108").Value = "03/01/2023"
$ws.Range("C108").Value = "15:00"
$ws.Range("D108").Value = "Levadeiakos"
$ws.Range("E108").Value = "Panathinaikos"
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = "A"
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 1
$ws.Range("K108").Value = "A"
$ws.Range("L108").Value = 14
$ws.Range("M108").Value = 8
$ws.Range("N108").Value = 4
$ws.Range("O108").Value = 3
$ws.Range("P108").Value = 17
$ws.Range("Q108").Value = 12
$ws.Range("R108").Value = 4
$ws.Range("S108").Value = 5
$ws.Range("T108").Value = 1
$ws.Range("U108").Value = 0
$ws.Range("V108").Value = 0
$ws.Range("W108").Value = 0
$ws.Range("X108").Value = 12
$ws.Range("Y108").Value = 4.5
$ws.Range("Z108").Value = 1.33
$ws.Range("AA108").Value = 2.15
$ws.Range("AB108").Value = 1.67

# Row 109: Giannina vs AEK (03/01/2023)
$ws.Range("A109").Value = "G1"
$ws.Range("B109").Value = "03/01/2023"
$ws.Range("C109").Value = "18:00"
$ws.Range("D109").Value = "Giannina"
$ws.Range("E109").Value = "AEK"
$ws.Range("F109").Value = 2
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = "H"
$ws.Range("I109").Value = 2
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = "H"
$ws.Range("L109").Value = 5
$ws.Range("M109").Value = 21
$ws.Range("N109").Value = 2
$ws.Range("O109").Value = 5
$ws.Range("P109").Value = 10
$ws.Range("Q109").Value = 15
$ws.Range("R109").Value = 0
$ws.Range("S109").Value = 6
$ws.Range("T109").Value = 3
$ws.Range("U109").Value = 1
$ws.Range("V109").Value = 0
$ws.Range("W109").Value = 0
$ws.Range("X109").Value = 11
$ws.Range("Y109").Value = 4.5
$ws.Range("Z109").Value = 1.33
$ws.Range("AA109").Value = 1.95
$ws.Range("AB109").Value = 1.9

# Row 110: Ionikos vs Olympiakos (03/01/2023)
$ws.Range("A110").Value = "G1"
$ws.Range("B110").Value = "03/01/2023"
$ws.Range("C110").Value = "19:30"
$ws.Range("D110").Value = "Ionikos"
$ws.Range("E110").Value = "Olympiakos"
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 2
$ws.Range("H110").Value = "A"
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 2
$ws.Range("K110").Value = "A"
$ws.Range("L110").Value = 2
$ws.Range("M110").Value = 15
$ws.Range("N110").Value = 1
$ws.Range("O110").Value = 5
$ws.Range("P110").Value = 11
$ws.Range("Q110").Value = 16
$ws.Range("R110").Value = 1
$ws.Range("S110").Value = 8
$ws.Range("T110").Value = 1
$ws.Range("U110").Value = 2
$ws.Range("V110").Value = 0
$ws.Range("W110").Value = 0
$ws.Range("X110").Value = 13
$ws.Range("Y110").Value = 5.25
$ws.Range("Z110").Value = 1.25
$ws.Range("AA110").Value = 1.8
$ws.Range("AB110").Value = 2

# Row 111: Asteras Tripolis vs Lamia (04/01/2023)
$ws.Range("A111").Value = "G1"
$ws.Range("B111").Value = "04/01/2023"
$ws.Range("C111").Value = "15:00"
$ws.Range("D111").Value = "Asteras Tripolis"
$ws.Range("E111").Value = "Lamia"
$ws.Range("F111").Value = 3
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = "H"
$ws.Range("I111").Value = 3
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = "H"
$ws.Range("L111").Value = 14
$ws.Range("M111").Value = 15
$ws.Range("N111").Value = 5
$ws.Range("O111").Value = 5
$ws.Range("P111").Value = 15
$ws.Range("Q111").Value = 14
$ws.Range("R111").Value = 3
$ws.Range("S111").Value = 2
$ws.Range("T111").Value = 3
$ws.Range("U111").Value = 3
$ws.Range("V111").Value = 0
$ws.Range("W111").Value = 0
$ws.Range("X111").Value = 1.75
$ws.Range("Y111").Value = 3.25
$ws.Range("Z111").Value = 5.5
$ws.Range("AA111").Value = 2.5
$ws.Range("AB111").Value = 1.5

# Row 112: Panetolikos vs Atromitos (04/01/2023)
$ws.Range("A112").Value = "G1"
$ws.Range("B112").Value = "04/01/2023"
$ws.Range("C112").Value = "16:00"
$ws.Range("D112").Value = "Panetolikos"
$ws.Range("E112").Value = "Atromitos"
$ws.Range("F112").Value = 2
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = "H"
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = "D"
$ws.Range("L112").Value = 15
$ws.Range("M112").Value = 6
$ws.Range("N112").Value = 5
$ws.Range("O112").Value = 1
$ws.Range("P112").Value = 14
$ws.Range("Q112").Value = 4
$ws.Range("R112").Value = 8
$ws.Range("S112").Value = 2
$ws.Range("T112").Value = 4
$ws.Range("U112").Value = 1
$ws.Range("V112").Value = 0
$ws.Range("W112").Value = 1
$ws.Range("X112").Value = 2.25
$ws.Range("Y112").Value = 3.2
$ws.Range("Z112").Value = 3.4
$ws.Range("AA112").Value = 2.25
$ws.Range("AB112").Value = 1.62

# Row 113: PAOK vs Aris (04/01/2023)
$ws.Range("A113").Value = "G1"
$ws.Range("B113").Value = "04/01/2023"
$ws.Range("C113").Value = "18:00"
$ws.Range("D113").Value = "PAOK"
$ws.Range("E113").Value = "Aris"
$ws.Range("F113").Value = 1
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = "H"
$ws.Range("I113").Value = 1
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = "H"
$ws.Range("L113").Value = 13
$ws.Range("M113").Value = 6
$ws.Range("N113").Value = 5
$ws.Range("O113").Value = 0
$ws.Range("P113").Value = 13
$ws.Range("Q113").Value = 20
$ws.Range("R113").Value = 9
$ws.Range("S113").Value = 2
$ws.Range("T113").Value = 4
$ws.Range("U113").Value = 1
$ws.Range("V113").Value = 0
$ws.Range("W113").Value = 0
$ws.Range("X113").Value = 1.91
$ws.Range("Y113").Value = 3.25
$ws.Range("Z113").Value = 4.5
$ws.Range("AA113").Value = 2.2
$ws.Range("AB113").Value = 1.65
